$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New header block P1:U1 - duplicate of the "parent info" headers (J1:O1).
$ws.Range("J1:O1").Copy($ws.Range("P1"))

# 2. New column block P:U on every data row holds a second guardian's info
#    (the "Ibu"/mother template that used to live in J3:O3). Fan it out to
#    rows 2-4 first, then fix up row 3's original J:O block with new data.
$ws.Range("J3:O3").Copy($ws.Range("P2"))
$ws.Range("J3:O3").Copy($ws.Range("P3"))
$ws.Range("J3:O3").Copy($ws.Range("P4"))

# 3. Row 3 (Annisa Putri) gets a new primary guardian (father instead of mother).
$ws.Range("J3").Value = "Budi Ahmad Jaya"
$ws.Range("K3").Value = "Ayah"
$ws.Range("L3").Value = "Jl. Merdeka No. 1, Jakarta"
$ws.Range("M3").Value = 81122334455
$ws.Range("N3").Value = "budi.aj@email.com"
$ws.Range("O3").Value = "Wiraswasta"

# 4. Personalize the second-guardian names per student.
$ws.Range("P3").Value = "Mariana"
$ws.Range("P4").Value = "Indah"

# 5. Turn the new e-mail address into a real hyperlink.
$ws.Hyperlinks.Add($ws.Range("N3"), "mailto:budi.aj@email.com", [Type]::Missing, [Type]::Missing, "budi.aj@email.com")

# 6. New column V: "siswa_baru_tingkat" (new-student grade level).
$ws.Range("V1").Value = "siswa_baru_tingkat"
$ws.Range("V1").Font.Name = "Arial"
$ws.Range("V1").Font.Size = 10
$ws.Range("V1").Borders.Item(7).Weight = -4138
$ws.Range("V1").Borders.Item(7).Color = 13421772
$ws.Range("V1").Borders.Item(10).Weight = -4138
$ws.Range("V1").Borders.Item(10).Color = 13421772

$ws.Range("V2").Value = 7
$ws.Range("V3").Value = 7
$ws.Range("V4").Value = 8

# 5. Update the view state (selection / scroll position) to match the saved file.
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("U7").Select()
